# Auto-generated: updates raw market-price values (columns H-N)
# across multiple sheets to match the scheduled-runner refresh.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 8087.9473
$ws.Range("I86").Value = 7996.7
$ws.Range("J86").Value = 8189.3335
$ws.Range("K86").Value = 7996.7
$ws.Range("L86").Value = 8189.3335
$ws.Range("M86").Value = -6873.7
$ws.Range("N86").Value = -10435.3335
$ws.Range("H89").Value = 8087.9473
$ws.Range("I89").Value = 7996.7
$ws.Range("J89").Value = 8189.3335
$ws.Range("K89").Value = 39983.5
$ws.Range("L89").Value = 40946.6675
$ws.Range("M89").Value = -34367.5
$ws.Range("N89").Value = -52178.6675
$ws.Range("H113").Value = 78000
$ws.Range("I113").Value = 78000
$ws.Range("J113").Value = 0
$ws.Range("K113").Value = 78000
$ws.Range("L113").Value = 0
$ws.Range("M113").Value = -74746
$ws.Range("N113").ClearContents()
$ws.Range("H137").Value = 14400.529
$ws.Range("J137").Value = 32589.143
$ws.Range("L137").Value = 97767.429
$ws.Range("N137").Value = -102867.429

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 5583.905
$ws.Range("I2").Value = 4757.636
$ws.Range("J2").Value = 6492.8
$ws.Range("K2").Value = 4757.636
$ws.Range("L2").Value = 6492.8
$ws.Range("M2").Value = -4644.636
$ws.Range("N2").Value = -6718.8
$ws.Range("H45").Value = 3170.0908
$ws.Range("I45").Value = 2743.5715
$ws.Range("K45").Value = 2743.5715
$ws.Range("M45").Value = -2366.5715
$ws.Range("H55").Value = 48666.668
$ws.Range("J55").Value = 48666.668
$ws.Range("L55").Value = 48666.668
$ws.Range("N55").Value = -49296.668
$ws.Range("H61").Value = 1868497
$ws.Range("I61").Value = 4571.6665
$ws.Range("J61").Value = 3965413
$ws.Range("K61").Value = 4571.6665
$ws.Range("L61").Value = 3965413
$ws.Range("M61").Value = -4359.6665
$ws.Range("N61").Value = -3965837
$ws.Range("H74").Value = 20641.846
$ws.Range("I74").Value = 1693.4667
$ws.Range("K74").Value = 1693.4667
$ws.Range("M74").Value = -819.4666999999999
$ws.Range("H77").Value = 20641.846
$ws.Range("I77").Value = 1693.4667
$ws.Range("K77").Value = 8467.333500000001
$ws.Range("M77").Value = -4099.333500000001
$ws.Range("H97").Value = 1217.5714
$ws.Range("I97").Value = 905
$ws.Range("J97").Value = 1634.3334
$ws.Range("K97").Value = 905
$ws.Range("L97").Value = 1634.3334
$ws.Range("M97").Value = -409
$ws.Range("N97").Value = -2626.3334
$ws.Range("H116").Value = 5583.905
$ws.Range("I116").Value = 4757.636
$ws.Range("J116").Value = 6492.8
$ws.Range("K116").Value = 4757.636
$ws.Range("L116").Value = 6492.8
$ws.Range("M116").Value = -2463.636
$ws.Range("N116").Value = -11080.8
$ws.Range("H122").Value = 3629.7058
$ws.Range("I122").Value = 1808.75
$ws.Range("K122").Value = 5426.25
$ws.Range("M122").Value = -2976.25
$ws.Range("H136").Value = 1868497
$ws.Range("I136").Value = 4571.6665
$ws.Range("J136").Value = 3965413
$ws.Range("K136").Value = 13714.9995
$ws.Range("L136").Value = 11896239
$ws.Range("M136").Value = -11164.9995
$ws.Range("N136").Value = -11901339

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 5583.905
$ws.Range("I3").Value = 4757.636
$ws.Range("J3").Value = 6492.8
$ws.Range("K3").Value = 4757.636
$ws.Range("L3").Value = 6492.8
$ws.Range("M3").Value = -4643.636
$ws.Range("N3").Value = -6720.8
$ws.Range("H76").Value = 0
$ws.Range("J76").Value = 0
$ws.Range("L76").Value = 0
$ws.Range("N76").ClearContents()
$ws.Range("H79").Value = 0
$ws.Range("J79").Value = 0
$ws.Range("L79").Value = 0
$ws.Range("N79").ClearContents()
$ws.Range("H99").Value = 19313
$ws.Range("I99").Value = 18726.924
$ws.Range("K99").Value = 18726.924
$ws.Range("M99").Value = -17228.924

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 10725.096
$ws.Range("I31").Value = 1107.7333
$ws.Range("J31").Value = 34768.5
$ws.Range("K31").Value = 1107.7333
$ws.Range("L31").Value = 34768.5
$ws.Range("M31").Value = -812.7333000000001
$ws.Range("N31").Value = -35358.5
$ws.Range("H34").Value = 10725.096
$ws.Range("I34").Value = 1107.7333
$ws.Range("J34").Value = 34768.5
$ws.Range("K34").Value = 1107.7333
$ws.Range("L34").Value = 34768.5
$ws.Range("M34").Value = -905.7333000000001
$ws.Range("N34").Value = -35172.5
$ws.Range("H86").Value = 17958.4
$ws.Range("J86").Value = 9998.5
$ws.Range("L86").Value = 9998.5
$ws.Range("N86").Value = -12244.5
$ws.Range("H89").Value = 17958.4
$ws.Range("J89").Value = 9998.5
$ws.Range("L89").Value = 49992.5
$ws.Range("N89").Value = -61224.5
$ws.Range("H122").Value = 2809.3333
$ws.Range("I122").Value = 1392.4445
$ws.Range("J122").Value = 3872
$ws.Range("K122").Value = 4177.333500000001
$ws.Range("L122").Value = 11616
$ws.Range("M122").Value = -1727.333500000001
$ws.Range("N122").Value = -16516
$ws.Range("H134").Value = 23260866
$ws.Range("I134").Value = 2016.4333
$ws.Range("J134").Value = 76935140
$ws.Range("K134").Value = 6049.2999
$ws.Range("L134").Value = 230805420
$ws.Range("M134").Value = -3514.2999
$ws.Range("N134").Value = -230810490

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H9").Value = 675331.25
$ws.Range("J9").Value = 575
$ws.Range("L9").Value = 1725
$ws.Range("N9").Value = -2173
$ws.Range("H68").Value = 996.25
$ws.Range("I68").Value = 799
$ws.Range("J68").Value = 1024.4286
$ws.Range("K68").Value = 2397
$ws.Range("L68").Value = 3073.2858
$ws.Range("M68").Value = -1586
$ws.Range("N68").Value = -4695.2858
$ws.Range("H71").Value = 996.25
$ws.Range("I71").Value = 799
$ws.Range("J71").Value = 1024.4286
$ws.Range("K71").Value = 7191
$ws.Range("L71").Value = 9219.857399999999
$ws.Range("M71").Value = -3135
$ws.Range("N71").Value = -17331.8574
$ws.Range("H107").Value = 2050.182
$ws.Range("I107").Value = 920.6
$ws.Range("J107").Value = 2991.5
$ws.Range("K107").Value = 2761.8
$ws.Range("L107").Value = 8974.5
$ws.Range("M107").Value = -841.8000000000002
$ws.Range("N107").Value = -12814.5
$ws.Range("H113").Value = 6750.5713
$ws.Range("I113").Value = 14288.556
$ws.Range("J113").Value = 1097.0834
$ws.Range("K113").Value = 42865.66800000001
$ws.Range("L113").Value = 3291.2502
$ws.Range("M113").Value = -40695.66800000001
$ws.Range("N113").Value = -7631.2502
$ws.Range("H121").Value = 1089.125
$ws.Range("I121").Value = 616
$ws.Range("K121").Value = 1848
$ws.Range("M121").Value = -538
$ws.Range("H131").Value = 3216.5715
$ws.Range("J131").Value = 2693.0544
$ws.Range("L131").Value = 8079.1632
$ws.Range("N131").Value = -18159.1632
$ws.Range("H139").Value = 9044.034
$ws.Range("I139").Value = 12559.77
$ws.Range("K139").Value = 37679.31
$ws.Range("M139").Value = -32539.31

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 2127.9473
$ws.Range("I97").Value = 2061.0833
$ws.Range("J97").Value = 2242.5715
$ws.Range("K97").Value = 2061.0833
$ws.Range("L97").Value = 2242.5715
$ws.Range("M97").Value = -1565.0833
$ws.Range("N97").Value = -3234.5715
$ws.Range("H123").Value = 54598.816
$ws.Range("J123").Value = 54598.816
$ws.Range("L123").Value = 54598.816
$ws.Range("N123").Value = -59498.816
$ws.Range("H126").Value = 11763.5
$ws.Range("I126").Value = 50000
$ws.Range("J126").Value = 7515
$ws.Range("K126").Value = 150000
$ws.Range("L126").Value = 22545
$ws.Range("M126").Value = -147530
$ws.Range("N126").Value = -27485
$ws.Range("H132").Value = 1913240.1
$ws.Range("I132").Value = 9999
$ws.Range("K132").Value = 29997
$ws.Range("M132").Value = -27467

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 6932.5557
$ws.Range("I7").Value = 7367.364
$ws.Range("J7").Value = 6249.2856
$ws.Range("K7").Value = 7367.364
$ws.Range("L7").Value = 6249.2856
$ws.Range("M7").Value = -7255.364
$ws.Range("N7").Value = -6473.2856
$ws.Range("H82").Value = 1966.0769
$ws.Range("I82").Value = 1862
$ws.Range("K82").Value = 1862
$ws.Range("M82").Value = -1501
$ws.Range("H85").Value = 1966.0769
$ws.Range("I85").Value = 1862
$ws.Range("K85").Value = 1862
$ws.Range("M85").Value = -614
$ws.Range("H100").Value = 2556.6
$ws.Range("I100").Value = 2309.8572
$ws.Range("J100").Value = 3132.3333
$ws.Range("K100").Value = 2309.8572
$ws.Range("L100").Value = 3132.3333
$ws.Range("M100").Value = -1768.8572
$ws.Range("N100").Value = -4214.3333
$ws.Range("H122").Value = 5513.448
$ws.Range("I122").Value = 4454.727
$ws.Range("J122").Value = 6160.4443
$ws.Range("K122").Value = 13364.181
$ws.Range("L122").Value = 18481.3329
$ws.Range("M122").Value = -10914.181
$ws.Range("N122").Value = -23381.3329
$ws.Range("H126").Value = 6932.5557
$ws.Range("I126").Value = 7367.364
$ws.Range("J126").Value = 6249.2856
$ws.Range("K126").Value = 22102.092
$ws.Range("L126").Value = 18747.8568
$ws.Range("M126").Value = -19632.092
$ws.Range("N126").Value = -23687.8568

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H42").Value = 0
$ws.Range("I42").Value = 0
$ws.Range("K42").Value = 0
$ws.Range("M42").ClearContents()
$ws.Range("H96").Value = 1499.1
$ws.Range("J96").Value = 1455.6666
$ws.Range("L96").Value = 1455.6666
$ws.Range("N96").Value = -4201.6666
$ws.Range("H100").Value = 807.1429000000001
$ws.Range("I100").Value = 880
$ws.Range("J100").Value = 625
$ws.Range("K100").Value = 1760
$ws.Range("L100").Value = 1250
$ws.Range("M100").Value = -1219
$ws.Range("N100").Value = -2332
$ws.Range("H107").Value = 3858.1667
$ws.Range("J107").Value = 900
$ws.Range("L107").Value = 2700
$ws.Range("N107").Value = -6540
$ws.Range("H127").Value = 24666.666
$ws.Range("J127").Value = 24666.666
$ws.Range("L127").Value = 24666.666
$ws.Range("N127").Value = -34586.666

